# docs: ajout d'avertissement étudiant et suppression mentions confidentielles
# Update the "date" column (C) for rows 2-21 on the active sheet: each
# timestamp is shifted forward by 15 hours 16 minutes (new run timestamps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "2026-02-18 09:29"
$ws.Range("C3").Value  = "2026-02-18 10:29"
$ws.Range("C4").Value  = "2026-02-18 11:29"
$ws.Range("C5").Value  = "2026-02-18 12:29"
$ws.Range("C6").Value  = "2026-02-18 13:29"
$ws.Range("C7").Value  = "2026-02-18 14:29"
$ws.Range("C8").Value  = "2026-02-18 15:29"
$ws.Range("C9").Value  = "2026-02-18 16:29"
$ws.Range("C10").Value = "2026-02-18 17:29"
$ws.Range("C11").Value = "2026-02-18 18:29"
$ws.Range("C12").Value = "2026-02-18 19:29"
$ws.Range("C13").Value = "2026-02-18 20:29"
$ws.Range("C14").Value = "2026-02-18 21:29"
$ws.Range("C15").Value = "2026-02-18 22:29"
$ws.Range("C16").Value = "2026-02-18 23:29"
$ws.Range("C17").Value = "2026-02-19 00:29"
$ws.Range("C18").Value = "2026-02-19 01:29"
$ws.Range("C19").Value = "2026-02-19 02:29"
$ws.Range("C20").Value = "2026-02-19 03:29"
$ws.Range("C21").Value = "2026-02-19 04:29"
